$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.610.95'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.460.53'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.78'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.21'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('E9').Value = '  +2.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.41'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.839.40'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.81'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.458.10'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.615.09'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0940'
$ws.Range('E20').Value = '  +1.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.58'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.31'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.86'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.91'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.27'
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.70'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.14'
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.84'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.46'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0759'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.42'
$ws.Range('E36').Value = '  -4.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.89'
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.103'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('E41').Value = '  -4.31%  '
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.975.15'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.86'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.00'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.696.98'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.88'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '66.85'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.27'
$ws.Range('E51').Value = '  +3.13%  '
